# Add Random forest algorithm predictions: insert two new columns (D, E)
# for "Predicted Eg" (Random Forest) values, shifting the existing
# "Predicted Eg_SVR" (D) and "Predicted Eg_RF" (E) columns to F and G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns before D. This shifts the existing D:E
#    (SVR / RF-label columns) to F:G, carrying over their exact widths,
#    per-cell styles and values intact.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) The insert leaves placeholder (empty, but styled) cells behind in the
#    new D:E columns for every row that had data in old D:E. Fully clear
#    D1:E35 so the new columns start out completely blank (no leftover
#    style), matching a genuinely "new" column.
$ws.Range("D1:E35").Clear()

# 3) Re-create the D1/E1 header cells ("Predicted Eg") using the same
#    look (bold / bordered / centered) as the neighbouring header C1.
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Predicted Eg"
$ws.Range("E1").Value = "Predicted Eg"

$ws.Range("D2").Value = 3.09
$ws.Range("E2").Value = 3.09
$ws.Range("D3").Value = 3.21
$ws.Range("E3").Value = 3.18
$ws.Range("D4").Value = 3.16
$ws.Range("E4").Value = 3.25
$ws.Range("D5").Value = 3.18
$ws.Range("E5").Value = 3.2
$ws.Range("D6").Value = 3.12
$ws.Range("E6").Value = 3.19
$ws.Range("D7").Value = 3.24
$ws.Range("E7").Value = 3.21
$ws.Range("D8").Value = 3.19
$ws.Range("E8").Value = 3.2
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 3.01
$ws.Range("D10").Value = 3.12
$ws.Range("E10").Value = 3.06
$ws.Range("D11").Value = 3.12
$ws.Range("E11").Value = 3.12
$ws.Range("D12").Value = 3.15
$ws.Range("E12").Value = 3.1
$ws.Range("D13").Value = 3.08
$ws.Range("E13").Value = 3.16
$ws.Range("D14").Value = 3.04
$ws.Range("E14").Value = 3.12
$ws.Range("D15").Value = 3.1
$ws.Range("E15").Value = 3.09
$ws.Range("D16").Value = 3.12
$ws.Range("E16").Value = 3.02
$ws.Range("D17").Value = 3.07
$ws.Range("E17").Value = 3.05
$ws.Range("D18").Value = 3.02
$ws.Range("E18").Value = 3.06
$ws.Range("D19").Value = 3.09
$ws.Range("E19").Value = 3.05
$ws.Range("D20").Value = 1.84
$ws.Range("E20").Value = 1.95
$ws.Range("D21").Value = 2.2799999999999998
$ws.Range("E21").Value = 2.69
$ws.Range("D22").Value = 2.2000000000000002
$ws.Range("E22").Value = 2.88
$ws.Range("D23").Value = 2.2400000000000002
$ws.Range("E23").Value = 2.86
$ws.Range("D24").Value = 2.2200000000000002
$ws.Range("E24").Value = 2.87
$ws.Range("D25").Value = 2.27
$ws.Range("E25").Value = 2.79
$ws.Range("D26").Value = 6.63
$ws.Range("E26").Value = 5.62
$ws.Range("D27").Value = 3.49
$ws.Range("E27").Value = 3.39
$ws.Range("D28").Value = 3.75
$ws.Range("E28").Value = 3.7
$ws.Range("D29").Value = 2.77
$ws.Range("E29").Value = 2.35
$ws.Range("D30").Value = 3.66
$ws.Range("E30").Value = 3.58
$ws.Range("D31").Value = 8.51
$ws.Range("E31").Value = 8.52
$ws.Range("D32").Value = 2.57
$ws.Range("E32").Value = 2.54
$ws.Range("D33").Value = 2.61
$ws.Range("E33").Value = 2.5499999999999998
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = 3
$ws.Range("D35").Value = 3.77
$ws.Range("E35").Value = 3.19

# 4) Update the active selection to match the edited workbook.
$ws.Range("J12").Select()

Write-Host "Random forest columns added"
